$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.097.71'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '2.689.09'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'559.22"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.45%  '
$ws.Range('D6').Value = "'159.33"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.73%  '
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('E10').Value = '  -2.12%  '
$ws.Range('E11').Value = '  -3.32%  '
$ws.Range('E12').Value = '  -5.94%  '
$ws.Range('D13').Value = '3.164.57'
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('E14').Value = '  -1.82%  '
$ws.Range('D15').Value = '62.999.76'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('E16').Value = '  -1.44%  '
$ws.Range('D17').Value = '2.689.87'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').Value = "'4.64"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.39%  '
$ws.Range('D20').Value = "'346.75"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('D21').Value = "'6.33"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.08%  '
$ws.Range('D22').Value = "'0.999"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('D24').Value = "'63.46"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('E25').Value = '  -0.71%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = "'8.27"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.85%  '
$ws.Range('D28').Value = "'1.45"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.41%  '
$ws.Range('D29').Value = '0.0₃0864'
$ws.Range('E29').Value = '  -4.90%  '
$ws.Range('D30').Value = "'7.28"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').Value = "'164.91"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.58%  '
$ws.Range('D33').Value = "'4.96"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.17%  '
$ws.Range('E34').Value = '  +1.33%  '
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = "'19.57"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('D38').Value = "'360.98"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.31%  '
$ws.Range('D39').Value = "'6.47"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.77%  '
$ws.Range('E40').Value = '  -2.35%  '
$ws.Range('D41').Value = "'4.03"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('D42').Value = "'38.55"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').Value = "'21.19"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.95%  '
$ws.Range('D44').Value = "'20.50"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').Value = "'0.622"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = "'0.0566"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.08%  '
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = "'11.05"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = "'130.06"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.14%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = "'0.0975"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.71%  '
